$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 389.0909
$ws.Range("I28").Value = 406.75
$ws.Range("J28").Value = 342
$ws.Range("K28").Value = 406.75
$ws.Range("L28").Value = 342
$ws.Range("M28").Value = 78.25
$ws.Range("N28").Value = -1312
$ws.Range("H112").Value = 2050.5
$ws.Range("J112").Value = 2185.8235
$ws.Range("L112").Value = 6557.470499999999
$ws.Range("N112").Value = -8773.470499999999
$ws.Range("H129").Value = 977.7755
$ws.Range("I129").Value = 490.07693
$ws.Range("K129").Value = 1470.23079
$ws.Range("M129").Value = 3529.76921
$ws.Range("H137").Value = 1433.7548
$ws.Range("I137").Value = 1151.1464
$ws.Range("J137").Value = 2399.3333
$ws.Range("K137").Value = 3453.4392
$ws.Range("L137").Value = 7197.999899999999
$ws.Range("M137").Value = -903.4392000000003
$ws.Range("N137").Value = -12297.9999
$ws.Range("H139").Value = 73485.71000000001
$ws.Range("J139").Value = 73485.71000000001
$ws.Range("L139").Value = 73485.71000000001
$ws.Range("N139").Value = -83765.71000000001
$ws.Range("H140").Value = 97220.625
$ws.Range("J140").Value = 97465.336
$ws.Range("L140").Value = 97465.336
$ws.Range("N140").Value = -107825.336

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13187.567
$ws.Range("I32").Value = 11871.131
$ws.Range("K32").Value = 11871.131
$ws.Range("M32").Value = -11584.131
$ws.Range("H122").Value = 3648.6365
$ws.Range("I122").Value = 3648.6365
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10945.9095
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -8495.9095
$ws.Range("N122").ClearContents()
$ws.Range("H131").Value = 25000
$ws.Range("J131").Value = 25000
$ws.Range("L131").Value = 25000
$ws.Range("N131").Value = -35080

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 23841.75
$ws.Range("J109").Value = 23841.75
$ws.Range("L109").Value = 23841.75
$ws.Range("N109").Value = -26615.75
$ws.Range("H140").Value = 49010.625
$ws.Range("J140").Value = 49010.625
$ws.Range("L140").Value = 49010.625
$ws.Range("N140").Value = -59370.625

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2630.6545
$ws.Range("I31").Value = 1415.3939
$ws.Range("J31").Value = 4453.5454
$ws.Range("K31").Value = 1415.3939
$ws.Range("L31").Value = 4453.5454
$ws.Range("M31").Value = -1120.3939
$ws.Range("N31").Value = -5043.5454
$ws.Range("H34").Value = 2630.6545
$ws.Range("I34").Value = 1415.3939
$ws.Range("J34").Value = 4453.5454
$ws.Range("K34").Value = 1415.3939
$ws.Range("L34").Value = 4453.5454
$ws.Range("M34").Value = -1213.3939
$ws.Range("N34").Value = -4857.5454
$ws.Range("H39").Value = 14887.75
$ws.Range("I39").Value = 12728.857
$ws.Range("K39").Value = 12728.857
$ws.Range("M39").Value = -12337.857
$ws.Range("H49").Value = 14887.75
$ws.Range("I49").Value = 12728.857
$ws.Range("K49").Value = 12728.857
$ws.Range("M49").Value = -12546.857
$ws.Range("H58").Value = 1428.7407
$ws.Range("I58").Value = 1577.7646
$ws.Range("J58").Value = 1175.4
$ws.Range("K58").Value = 1577.7646
$ws.Range("L58").Value = 1175.4
$ws.Range("M58").Value = -1374.7646
$ws.Range("N58").Value = -1581.4
$ws.Range("H136").Value = 1428.7407
$ws.Range("I136").Value = 1577.7646
$ws.Range("J136").Value = 1175.4
$ws.Range("K136").Value = 4733.293799999999
$ws.Range("L136").Value = 3526.2
$ws.Range("M136").Value = -2183.293799999999
$ws.Range("N136").Value = -8626.200000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 11355.4
$ws.Range("I70").Value = 13194.25
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 39582.75
$ws.Range("L70").Value = 12000
$ws.Range("M70").Value = -39267.75
$ws.Range("N70").Value = -12630
$ws.Range("H73").Value = 11355.4
$ws.Range("I73").Value = 13194.25
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 39582.75
$ws.Range("L73").Value = 12000
$ws.Range("M73").Value = -38490.75
$ws.Range("N73").Value = -14184
$ws.Range("J113").Value = 6655.5
$ws.Range("L113").Value = 19966.5
$ws.Range("N113").Value = -24306.5
$ws.Range("H131").Value = 10640724
$ws.Range("I131").Value = 9505.833000000001
$ws.Range("J131").Value = 12196512
$ws.Range("K131").Value = 28517.499
$ws.Range("L131").Value = 36589536
$ws.Range("M131").Value = -23477.499
$ws.Range("N131").Value = -36599616
$ws.Range("H132").Value = 1942.0555
$ws.Range("I132").Value = 1229
$ws.Range("J132").Value = 2395.818
$ws.Range("K132").Value = 11061
$ws.Range("L132").Value = 21562.362
$ws.Range("M132").Value = -8531
$ws.Range("N132").Value = -26622.362

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 22700
$ws.Range("J18").Value = 13600
$ws.Range("L18").Value = 13600
$ws.Range("N18").Value = -14186
$ws.Range("H51").Value = 28642.572
$ws.Range("J51").Value = 28642.572
$ws.Range("L51").Value = 28642.572
$ws.Range("N51").Value = -29660.572
$ws.Range("H70").Value = 5740.433
$ws.Range("I70").Value = 4978.385
$ws.Range("J70").Value = 6323.1763
$ws.Range("K70").Value = 4978.385
$ws.Range("L70").Value = 6323.1763
$ws.Range("M70").Value = -4708.385
$ws.Range("N70").Value = -6863.1763
$ws.Range("H73").Value = 5740.433
$ws.Range("I73").Value = 4978.385
$ws.Range("J73").Value = 6323.1763
$ws.Range("K73").Value = 4978.385
$ws.Range("L73").Value = 6323.1763
$ws.Range("M73").Value = -4042.385
$ws.Range("N73").Value = -8195.176299999999
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H122").Value = 3967.7827
$ws.Range("I122").Value = 3806.7742
$ws.Range("J122").Value = 4300.533
$ws.Range("K122").Value = 11420.3226
$ws.Range("L122").Value = 12901.599
$ws.Range("M122").Value = -8970.3226
$ws.Range("N122").Value = -17801.599
$ws.Range("H140").Value = 50612
$ws.Range("J140").Value = 50612
$ws.Range("L140").Value = 50612
$ws.Range("N140").Value = -60972

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H131").Value = 21057.555
$ws.Range("I131").Value = 8296
$ws.Range("J131").Value = 24703.715
$ws.Range("K131").Value = 8296
$ws.Range("L131").Value = 24703.715
$ws.Range("M131").Value = -3256
$ws.Range("N131").Value = -34783.715
$ws.Range("H136").Value = 3267.125
$ws.Range("I136").Value = 3239.9272
$ws.Range("J136").Value = 3433.3333
$ws.Range("K136").Value = 9719.7816
$ws.Range("L136").Value = 10299.9999
$ws.Range("M136").Value = -7169.7816
$ws.Range("N136").Value = -15399.9999
$ws.Range("H138").Value = 88960
$ws.Range("J138").Value = 88960
$ws.Range("L138").Value = 88960
$ws.Range("N138").Value = -99240

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3420
$ws.Range("I62").Value = 3300
$ws.Range("J62").Value = 3500
$ws.Range("K62").Value = 3300
$ws.Range("L62").Value = 3500
$ws.Range("M62").Value = -2676
$ws.Range("N62").Value = -4748
$ws.Range("H65").Value = 3420
$ws.Range("I65").Value = 3300
$ws.Range("J65").Value = 3500
$ws.Range("K65").Value = 16500
$ws.Range("L65").Value = 17500
$ws.Range("M65").Value = -13380
$ws.Range("N65").Value = -23740
$ws.Range("H123").Value = 24307.166
$ws.Range("J123").Value = 24307.166
$ws.Range("L123").Value = 24307.166
$ws.Range("N123").Value = -34107.166
$ws.Range("H136").Value = 1670.7916
$ws.Range("I136").Value = 1551.2565
$ws.Range("J136").Value = 2188.7778
$ws.Range("K136").Value = 4653.7695
$ws.Range("L136").Value = 6566.3334
$ws.Range("M136").Value = -2103.7695
$ws.Range("N136").Value = -11666.3334
